$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_5a_Indikatoren")

# --- Column width adjustments (col J/10, L/12, M/13) ---
$ws.Columns(10).ColumnWidth = 51.857142857142854
$ws.Columns(12).ColumnWidth = 48
$ws.Columns(13).ColumnWidth = 45

# --- Cell text updates (Ziel kurz De / Ziel kurz En, plus a few De/En fixes) ---
$ws.Range('M2').Value = 'Keep the proportion considerably below the EU-level by 2030'
$ws.Range('M3').Value = 'Keep the proportion considerably below the EU-level by 2030'
$ws.Range('M4').Value = 'Reduction of the nitrogen surpluses to 70 kg/ha on an annual average between 2028 and 2032'
$ws.Range('L5').Value = 'Erhöhung des Anteils von ökologischeM Landbau auf 20 % bis 2030'
$ws.Range('M5').Value = 'Increase the proportion of organically farmed agricultural land to 20% by 2030'
$ws.Range('M6').Value = 'Funds disbursed for the application of the guidelines and recommendations of the UN Committee on World Food Security (CFS) to be increased appropriately by 2030'
$ws.Range('L7').Value = 'Senkung auf 100 Todesfälle je 100 000 Einwohnerinnen bis 2030'
$ws.Range('M7').Value = 'To be reduced to 100 deaths per 100,000 inhabitants by 2030'
$ws.Range('L8').Value = 'Senkung auf 190 Todesfälle je 100 000 Einwohner bis 2030'
$ws.Range('M8').Value = 'To be reduced to 190 deaths per 100,000 inhabitants by 2030'
$ws.Range('M9').Value = 'Reduction to 7% by 2030'
$ws.Range('M10').Value = 'Reduction to 19% by 2030'
$ws.Range('M11').Value = 'Increase to be permanently halted'
$ws.Range('M12').Value = 'Increase to be permanently halted'
$ws.Range('M13').Value = 'Reduction to 55% by 2030'
$ws.Range('L14').Value = 'Erreichung des Feinstaub-Richtwerts von 20 µg/m³ bis 2030'
$ws.Range('M14').Value = 'Adherence to the particulate matter guideline value by 2030'
$ws.Range('M15').Value = 'Expenditure to be increased by 2030'
$ws.Range('L16').Value = 'Verringerung des Anteils auf 9,5 % bis 2030'
$ws.Range('M16').Value = 'Reduce the proportion to 9.5% by 2030'
$ws.Range('L17').Value = 'Steigerung des Anteils auf 55 % bis 2030'
$ws.Range('M17').Value = 'Increase the proportion to 55% by 2030'
$ws.Range('L18').Value = 'Anstieg auf 35 % bis 2030'
$ws.Range('M18').Value = 'Increase to 35% by 2030'
$ws.Range('L19').Value = 'Anstieg auf 60 % bis 2020 und auf 70 % bis 2030'
$ws.Range('M19').Value = 'Increase to 60% by 2020 and to 70% by 2030'
$ws.Range('L20').Value = 'Verringerung des Abstandes auf 10 % bis 2020, Beibehaltung bis 2030'
$ws.Range('M20').Value = 'Reduce the gap to 10% by 2020, maintained until 2030'
$ws.Range('L21').Value = '30 % Frauen in Aufsichtsräten bis 2030'
$ws.Range('M21').Value = '30% women on supervisory boards by 2030'
$ws.Range('L22').Value = 'Gleichberechtigte Teilhabe im öffentlichen Dienst bis 2025'
$ws.Range('M22').Value = 'Equal-opportunity participation in civil service by 2025'
$ws.Range('L23').Value = '65 % bis 2030'
$ws.Range('M23').Value = '65% by 2030'
$ws.Range('L24').Value = 'Steigerung um ein Drittel bis 2030'
$ws.Range('M24').Value = 'To be increased by a third by 2030'
$ws.Range('L25').Value = 'Einhaltung der gewässertypischen Orientierungswerte bis 2030'
$ws.Range('M25').Value = 'Not exceeding benchmark values by 2030'
$ws.Range('L26').Value = 'Einhaltung des Schwellenwertes bis 2030'
$ws.Range('M26').Value = 'Compliance with the threshold value by 2030'
$ws.Range('L28').Value = '6 Millionen Menschen pro Jahr bis 2030'
$ws.Range('M28').Value = '6 million people per year by 2030'
$ws.Range('L29').Value = '4 Millionen Menschen pro Jahr bis 2030'
$ws.Range('M29').Value = '4 million people per year by 2030'
$ws.Range('L30').Value = 'Steigerung um 2,1 % pro Jahr'
$ws.Range('M30').Value = 'Increase by 2.1% per year'
$ws.Range('L31').Value = 'Senkung um 20 % bis 2020, um 30 % bis 2030 und um 50 % bis 2050'
$ws.Range('M31').Value = 'Reduction by 20% by 2020, by 30% by 2030, and by 50% by 2050'
$ws.Range('L32').Value = 'Anstieg auf 18 % bis 2020, auf 30 % bis 2030, auf 45 % bis 2040 und auf 60 % bis 2050'
$ws.Range('M32').Value = 'Increase to 18% by 2020 and to 30% by 2030, to 45% by 2040 and to 60% by 2050'
$ws.Range('L33').Value = 'Anstieg auf 35 % bis 2020, 65 % bis 2030 und Treibhausgasneutralität bis 2050'
$ws.Range('M33').Value = 'Increase to 35% by 2020, 65% by 2030, and greenhouse gas neutrality by 2050'
$ws.Range('L34').Value = 'Beibehaltung des Trends der Jahre 2000 – 2010 bis 2030'
$ws.Range('M34').Value = 'Trend of the years 2000–2010 to be maintained until 2030'
$ws.Range('L35').Value = 'Unter 3 % des BIP'
$ws.Range('M35').Value = 'Less than 3% of GDP'
$ws.Range('L36').Value = 'Unter 0,5 % des BIP'
$ws.Range('M36').Value = 'Less than 0.5% of GDP'
$ws.Range('L37').Value = 'Unter 60 % des BIP'
$ws.Range('M37').Value = 'Less than 60% of GDP'
$ws.Range('L38').Value = 'Angemessene Entwicklung bis 2030'
$ws.Range('M38').Value = 'Appropriate development by 2030'
$ws.Range('L39').Value = 'Stetiges und angemessenes Wirtschaftswachstum'
$ws.Range('M39').Value = 'Steady and appropriate economic growth'
$ws.Range('L40').Value = 'Erhöhung auf 78 % bis 2030'
$ws.Range('M40').Value = 'Increase to 78% by 2030'
$ws.Range('L41').Value = 'Erhöhung auf 60 % bis 2030'
$ws.Range('M41').Value = 'Increase to 60% by 2030'
$ws.Range('L42').Value = 'Signifikante Steigerung bis 2030'
$ws.Range('M42').Value = 'Significantly increase by 2030'
$ws.Range('L43').Value = 'Jährlich mindestens 3,5 % des BIP bis 2025'
$ws.Range('M43').Value = 'At least 3.5% of GDP per year by 2025'
$ws.Range('L44').Value = 'Flächendeckend bis 2025'
$ws.Range('M44').Value = 'Full coverage by 2025'
$ws.Range('L45').Value = 'Erhöhung der Abschlussquote ausländischen Schulabgänger und Angleichung an die Quote deutscher Abgänger bis 2030'
$ws.Range('M45').Value = 'Increase the graduation rate of foreign school leavers and bring it into line with the rate of German school leavers by 2030'
$ws.Range('J46').Value = 'Gini-Koeffizient Einkommen nach Sozialtransfer bis 2030 unterhalb des EU-Wertes'
$ws.Range('K46').Value = 'Gini coefficient of income after social transfers to be below the EU figure by 2030'
$ws.Range('L46').Value = 'Bis 2030 unterhalb des EU-Wertes halten'
$ws.Range('M46').Value = 'To be below the EU figure by 2030'
$ws.Range('L47').Value = 'Senkung auf 30 ha pro Tag bis 2030'
$ws.Range('M47').Value = 'Reduction to under 30 ha per day by 2030'
$ws.Range('L48').Value = 'Verringerung'
$ws.Range('M48').Value = 'Reduce the loss'
$ws.Range('L49').Value = 'Keine Verringerung'
$ws.Range('M49').Value = 'No reduction'
$ws.Range('L50').Value = 'Senkung um 15 bis 20 % bis 2030'
$ws.Range('M50').Value = 'Reduction by 15–20% by 2030'
$ws.Range('L51').Value = 'Senkung um 15 bis 20 % bis 2030'
$ws.Range('M51').Value = 'Reduction by 15–20% by 2030'
$ws.Range('L52').Value = 'Verringerung der Reisezeit mit öffentlichen Verkehrsmitteln'
$ws.Range('M52').Value = 'Reduction of travel time by public transport'
$ws.Range('L53').Value = 'Senkung auf 13 % bis 2030'
$ws.Range('M53').Value = 'Reduce to 13% by 2030'
$ws.Range('L54').Value = 'Steigerung auf 50 Millionen bis 2030'
$ws.Range('M54').Value = 'Increase to 50 million by 2030'
$ws.Range('L55').Value = 'Steigerung auf 34 % bis 2030'
$ws.Range('M55').Value = 'Increase to 34% by 2030'
$ws.Range('L56').Value = 'Kontinuierliche Reduzierung '
$ws.Range('M56').Value = 'Steady reduction'
$ws.Range('L57').Value = 'Kontinuierliche Reduzierung '
$ws.Range('M57').Value = 'Steady reduction'
$ws.Range('L58').Value = 'Kontinuierliche Reduzierung '
$ws.Range('M58').Value = 'Steady reduction'
$ws.Range('L59').Value = '5 000 Standorte bis 2030'
$ws.Range('M59').Value = '5,000 locations by 2030'
$ws.Range('L60').Value = 'Steigerung auf 95 % bis 2020'
$ws.Range('M60').Value = 'Increase to 95% by 2020'
$ws.Range('L61').Value = 'Signifikante Senkung'
$ws.Range('M61').Value = 'Significantly reduce'
$ws.Range('L62').Value = 'Minderung um 40 % bis 2020, 55 % bis 2030 und Treibhausgasneutralität bis 2050'
$ws.Range('M62').Value = 'Reduce by 40% by 2020, 55% by 2030 and greenhouse gas neutrality by 2050'
$ws.Range('L63').Value = 'Verdopplung bis 2020'
$ws.Range('M63').Value = 'Double by 2020'
$ws.Range('L64').Value = 'Gesamtstickstoff in Zuflüssen unter 2,6 mg/l'
$ws.Range('M64').Value = 'Total nitrogen in the inflows below 2.6 mg/l'
$ws.Range('L65').Value = 'Gesamtstickstoff in Zuflüssen unter 2,8 mg/l'
$ws.Range('M65').Value = 'Total nitrogen in the inflows below 2.8 mg/l'
$ws.Range('L66').Value = 'Außschließlich nachhaltig bewirtschaftete Fischbestände bis 2020'
$ws.Range('M66').Value = 'Only sustainably managed fish stocks by 2020'
$ws.Range('L67').Value = 'Erreichen des Indexwertes 100 bis 2030'
$ws.Range('M67').Value = 'Reach the index value of 100 by 2030'
$ws.Range('L68').Value = 'Verringerung um 35 % bis 2030'
$ws.Range('M68').Value = 'Reduction by 35% by 2030'
$ws.Range('L69').Value = 'Steigerung der Zahlungen bis 2030'
$ws.Range('M69').Value = 'Increase payments by 2030'
$ws.Range('L70').Value = 'Steigerung der Zahlungen bis 2030'
$ws.Range('M70').Value = 'Payments to be increased by 2030'
$ws.Range('L71').Value = 'Rückgang von Straftaten auf unter 6 500 je 100 000 Einwohner/ -innen bis 2030'
$ws.Range('M71').Value = 'Reduce criminal offences  to less than 6,500 per 100,000 inhabitants by 2030'
$ws.Range('L72').Value = 'Mindestens 15 Projekte pro Jahr bis 2030'
$ws.Range('M72').Value = 'At least 15 projects per year by 2030'
$ws.Range('L73').Value = 'Verbesserung bis 2030'
$ws.Range('M73').Value = 'Improvement by 2030'
$ws.Range('L74').Value = 'Verbesserung bis 2030'
$ws.Range('M74').Value = 'Improvement by 2030'
$ws.Range('L75').Value = 'Steigerung auf 0,7 % bis 2030'
$ws.Range('M75').Value = 'Increase to 0.7% by 2030'
$ws.Range('L76').Value = 'Steigerung um 10 % von 2015 bis 2020, anschließend Verstetigung'
$ws.Range('M76').Value = 'Increase by 10% from 2015 to 2020, then stabilised'
$ws.Range('L77').Value = 'Steigerung um 100 % bis 2030'
$ws.Range('M77').Value = 'Increase by 100 % by 2030'
